$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
Write-Output $ws.Name
